$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stale "_GoBack" bookmark (Word drops this automatically
#    once a new edit is made elsewhere in the document).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Replace "Ben changing things up!" with the new sentence, keeping
#    the trailing space as its own separate run (as in the target
#    document) instead of merged into the sentence's run.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Ben changing things up!", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Katherine creates a black-box test.", 2)

# Insert a temporary carrier paragraph right after the sentence we just
# wrote, put a single space in it, then cut that space out, delete the
# now-empty carrier paragraph mark (merging it away without leaving any
# trace), and paste the space back in as its own run immediately after
# the sentence. This reproduces the two-run layout
#   <w:r><w:t>Katherine creates a black-box test.</w:t></w:r>
#   <w:r><w:t xml:space="preserve"> </w:t></w:r>
# without disturbing the untouched paragraphs that follow.
$target = $d.Content
$target.Find.Execute("Katherine creates a black-box test.")
$target.Collapse(0)
$target.InsertParagraphAfter()

$carrier = $d.Paragraphs(6)
$carrier.Range.Text = " "

$carrier = $d.Paragraphs(6)
$spaceRng = $d.Range($carrier.Range.Start, $carrier.Range.Start + 1)
$spaceRng.Cut()

$carrier = $d.Paragraphs(6)
$carrier.Range.Delete()

$target2 = $d.Content
$target2.Find.Execute("Katherine creates a black-box test.")
$target2.Collapse(0)
$target2.Paste()
